# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column (populated
# with the literal value "stock" on every data row), inserted right after
# the existing "total" column and before the existing "date" column. The
# columns that used to be H:J ("date", "legislator_name", "legislator_id")
# shift right to I:K.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H; this shifts the old H:J ("date", "legislator_name",
# "legislator_id") to I:K and the new column inherits formatting from its
# neighbours (matches the bold/bordered header style in row 1 and the plain
# data style in the rows below).
$ws.Columns.Item(8).Insert()

# Header for the newly inserted column.
$ws.Range("H1").Value = "property_category"

# Every stock row gets the same property_category value: "stock".
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
